$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (SIQ_PO1_DGC_CYRS_006): the Return Date (H7) and Answer (I7) had been
# filled in but the reviewer wants them cleared out again - the answer
# wasn't actually returned/accepted yet, so only the Expected Return Date
# (G7) stays populated. Restyle G7 to match the "expected date only" look
# used on F7 (no special left-align), and restyle the now-empty H7/I7 cells
# to the plain wrap-text look already used elsewhere in the sheet (e.g. H10).
$ws.Range("F7").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("H10").Copy()
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("I7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H7:I7").ClearContents()

# Row 8 (SIQ_PO1_DGC_CYRS_007): same treatment.
$ws.Range("F8").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("H10").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("I8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H8:I8").ClearContents()

# Row 9 (SIQ_PO1_DGC_CYRS_008): same treatment.
$ws.Range("F9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("H10").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("I9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H9:I9").ClearContents()

# Reflect where the author's cursor ended up after reviewing the sheet.
$ws.Range("A9").Select()
